$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cells($ws, $cellA, $cellB, $scratch) {
    $cellA.Copy($scratch)
    $cellB.Copy($cellA)
    $scratch.Copy($cellB)
}

function Swap-Rows($ws, $rowA, $rowB, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $scratch = $ws.Cells.Item(500, $col)
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        Swap-Cells $ws $cellA $cellB $scratch
        $scratch.Clear()
    }
}

# Swap row 5 (Jose Alvarado) and row 6 (CJ McCollum) data across columns B..K
# (column A, the "No." index, stays fixed per physical row).
Swap-Rows $ws 5 6 2 11

# Swap row 14 (Garrett Temple) and row 15 (Kira Lewis Jr.) data across columns B..K
Swap-Rows $ws 14 15 2 11
